# Fleet (flotta) operator-reassignment update.
#
# Two vehicles that were unassigned get a new operator, effective 2026-01-28,
# and that change is logged in the "Storico Passaggi" (change history) sheet,
# which keeps only the two most-recent change entries (rows 2 and 3 are
# overwritten in place).
#
#   GL590TH  (row 72 on "Stato Attuale"): DA ASSEGNARE (BENEGIAMO MALATTIA LUNGA) -> BRAGHINI MATTEO
#   GL599TH  (row 79 on "Stato Attuale"): DA ASSEGNARE                            -> MASCIARELLI MAURIZIO
#
# Note: the date columns in this workbook store dates as plain text
# ("YYYY-MM-DD" shared strings), not real Excel date serials. Assigning a
# date-look-alike string straight to .Value makes Excel auto-convert it to a
# real date (and stamps a new number-format style on the cell), so instead we
# write it as a text formula and paste-special the literal value back in,
# which keeps the cell a plain shared string with no style changes - exactly
# matching how the rest of the sheet stores its dates.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Stato Attuale")
$ws2 = $wb.Worksheets.Item("Storico Passaggi")

$newDate = "2026-01-28"
$xlPasteValues = -4163

function Set-TextDate($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

# --- "Stato Attuale": reassign the two pending vehicles ---

# GL590TH
$ws1.Range("B72").Value = "BRAGHINI MATTEO"
Set-TextDate $ws1.Range("C72") $newDate

# GL599TH
$ws1.Range("B79").Value = "MASCIARELLI MAURIZIO"
Set-TextDate $ws1.Range("C79") $newDate

# --- "Storico Passaggi": log both changes (rows 2-3 hold the latest entries) ---

$ws2.Range("A2").Value = "GL590TH"
$ws2.Range("B2").Value = "DA ASSEGNARE (BENEGIAMO MALATTIA LUNGA)"
$ws2.Range("C2").Value = "BRAGHINI MATTEO"
Set-TextDate $ws2.Range("D2") $newDate

$ws2.Range("A3").Value = "GL599TH"
$ws2.Range("B3").Value = "DA ASSEGNARE"
$ws2.Range("C3").Value = "MASCIARELLI MAURIZIO"
Set-TextDate $ws2.Range("D3") $newDate
